# This edit re-orders the data rows 9-16 of the "Artfynd" sheet: the
# observation record that used to live in one row moves to another row
# (a permutation of whole rows, not a change of individual field values).
#
# Mapping: for each destination row number, the value is the row number
# (in the *original* layout) whose full record should end up there.
#   9  <- 10
#   10 <- 15
#   11 <- 13
#   12 <- 11
#   13 <- 12
#   14 <- 16
#   15 <- 9
#   16 <- 14
#
# Because this is a permutation (not simple pairwise swaps), every
# affected row is first read into memory in full, and only then written
# back out in the new order - this avoids clobbering source data that is
# still needed for a later row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # column A
$lastCol  = 51  # column AY

$rowMap = @{
    9  = 10
    10 = 15
    11 = 13
    12 = 11
    13 = 12
    14 = 16
    15 = 9
    16 = 14
}

# 1) Snapshot every source row that participates in the permutation.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowData = @()
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $rowData += ,($ws.Cells.Item($srcRow, $col).Value2)
        }
        $snapshot[$srcRow] = $rowData
    }
}

# 2) Write each destination row from the snapshot of its source row.
#    Plain-text cells that look like ISO dates ("2023-08-15") are stored
#    in this sheet as literal text, not real date values. Re-assigning
#    such a string straight back through .Value2 would let Excel's
#    automatic data-type detection reinterpret it as a date serial, so
#    those are re-written with a leading apostrophe to force text, just
#    like typing '2023-08-15 into a cell would.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $val = $rowData[$col - $firstCol]
        if ($val -is [string] -and $val -match '^\d{4}-\d{2}-\d{2}$') {
            $ws.Cells.Item($destRow, $col).Value2 = "'" + $val
        } else {
            $ws.Cells.Item($destRow, $col).Value2 = $val
        }
    }
}
